$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "162003.10+552822.4" (row 5) entirely,
# shifting the rows below it (row 6) up by one.
$ws.Rows("5:5").Delete()

# Update the active cell selection to H8, as recorded after the edit.
$ws.Range("H8").Select()
